$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the 2020-11-10 refresh of the "volet 2" regional / legal-category
# breakdown: nombre_aides (col C) and montant_total (col D) for a handful of
# region/categorie-juridique rows.
#
# Both columns are stored as *text* in this workbook (General-formatted cells
# holding numeric-looking strings), so a plain `.Value = "441"` assignment
# would let Excel auto-coerce the cell to a real number, which is not what we
# want. To keep the cell's data type as text (matching every other cell in
# the sheet) we stage the new text in an untouched helper cell that has been
# forced to text via a quote-prefix, copy it, and paste-special just the
# values into the destination cell - that brings the text across without
# flipping the destination's type to numeric. The helper cell is fully
# cleared afterwards so it leaves no trace in the sheet.

$updates = @(
    @{Row=31;  C="441";  D="1300749.11"},
    @{Row=33;  C="819";  D="5354648.45"},
    @{Row=35;  C="548";  D="2974242.32"},
    @{Row=38;  C="591";  D="1547372.32"},
    @{Row=39;  C="293";  D="1590068.04"},
    @{Row=40;  C="275";  D="928520.72"},
    @{Row=50;  C="996";  D="6372799.81"},
    @{Row=58;  C="6950"; D="35597073.29"},
    @{Row=60;  C="6812"; D="29488918.34"},
    @{Row=62;  C="137";  D="692707.46"},
    @{Row=100; C="1376"; D="3483754.28"},
    @{Row=103; C="1563"; D="7975645.55"},
    @{Row=105; C="1509"; D="7035240.06"}
)

$helper = $ws.Range("ZZ1")

foreach ($u in $updates) {
    foreach ($col in @("C", "D")) {
        $newText = $u[$col]
        $target = $ws.Range($col + $u.Row)

        $helper.Value = "'" + $newText
        $helper.Copy()
        $target.PasteSpecial(-4163)  # xlPasteValues
    }
}

$helper.Clear()
$excel.CutCopyMode = 0

Write-Output "Updated nombre_aides / montant_total for rows: 31,33,35,38,39,40,50,58,60,62,100,103,105"
